$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 969.1
$ws.Range("J17").Value = 969.1
$ws.Range("L17").Value = 2907.3
$ws.Range("N17").Value = -3243.3

$ws.Range("H40").Value = 2003.2174
$ws.Range("I40").Value = 1380
$ws.Range("J40").Value = 2482.6155
$ws.Range("K40").Value = 1380
$ws.Range("L40").Value = 2482.6155
$ws.Range("M40").Value = -1205
$ws.Range("N40").Value = -2832.6155

$ws.Range("H70").Value = 1019.53656
$ws.Range("I70").Value = 1202.6296
$ws.Range("J70").Value = 666.4286
$ws.Range("K70").Value = 3607.8888
$ws.Range("L70").Value = 1999.2858
$ws.Range("M70").Value = -3337.8888
$ws.Range("N70").Value = -2539.2858

$ws.Range("H73").Value = 1019.53656
$ws.Range("I73").Value = 1202.6296
$ws.Range("J73").Value = 666.4286
$ws.Range("K73").Value = 3607.8888
$ws.Range("L73").Value = 1999.2858
$ws.Range("M73").Value = -2671.8888
$ws.Range("N73").Value = -3871.2858

$ws.Range("H86").Value = 101488.164
$ws.Range("I86").Value = 201250.5
$ws.Range("J86").Value = 1725.8334
$ws.Range("K86").Value = 201250.5
$ws.Range("L86").Value = 1725.8334
$ws.Range("M86").Value = -200127.5
$ws.Range("N86").Value = -3971.8334

$ws.Range("H89").Value = 101488.164
$ws.Range("I89").Value = 201250.5
$ws.Range("J89").Value = 1725.8334
$ws.Range("K89").Value = 1006252.5
$ws.Range("L89").Value = 8629.166999999999
$ws.Range("M89").Value = -1000636.5
$ws.Range("N89").Value = -19861.167

$ws.Range("H97").Value = 2600
$ws.Range("I97").Value = 2000
$ws.Range("J97").Value = 2800
$ws.Range("K97").Value = 6000
$ws.Range("L97").Value = 8400
$ws.Range("M97").Value = -5504
$ws.Range("N97").Value = -9392

$ws.Range("H99").Value = 633.0833
$ws.Range("I99").Value = 282.55554
$ws.Range("J99").Value = 1684.6666
$ws.Range("K99").Value = 847.66662
$ws.Range("L99").Value = 5053.9998
$ws.Range("M99").Value = 650.33338
$ws.Range("N99").Value = -8049.9998

$ws.Range("H101").Value = 2206.8572
$ws.Range("I101").Value = 1703.7142
$ws.Range("J101").Value = 2710
$ws.Range("K101").Value = 5111.142599999999
$ws.Range("L101").Value = 8130
$ws.Range("M101").Value = -3489.142599999999
$ws.Range("N101").Value = -11374

$ws.Range("H138").Value = 1763.2933
$ws.Range("I138").Value = 1330.4773
$ws.Range("J138").Value = 2377.6128
$ws.Range("K138").Value = 3991.4319
$ws.Range("L138").Value = 7132.8384
$ws.Range("M138").Value = 1148.5681
$ws.Range("N138").Value = -17412.8384

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2177.875
$ws.Range("I94").Value = 2145.077
$ws.Range("J94").Value = 2320
$ws.Range("K94").Value = 2145.077
$ws.Range("L94").Value = 2320
$ws.Range("M94").Value = -1694.077
$ws.Range("N94").Value = -3222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2828.7222
$ws.Range("I7").Value = 9148.362999999999
$ws.Range("J7").Value = 48.08
$ws.Range("K7").Value = 9148.362999999999
$ws.Range("L7").Value = 48.08
$ws.Range("M7").Value = -9035.362999999999
$ws.Range("N7").Value = -274.08

$ws.Range("H11").Value = 1100
$ws.Range("J11").Value = 1100
$ws.Range("L11").Value = 1100
$ws.Range("N11").Value = -1380

$ws.Range("H43").Value = 23225
$ws.Range("J43").Value = 23225
$ws.Range("L43").Value = 23225
$ws.Range("N43").Value = -23593

$ws.Range("H74").Value = 18304.666
$ws.Range("J74").Value = 18304.666
$ws.Range("L74").Value = 18304.666
$ws.Range("N74").Value = -20052.666

$ws.Range("H77").Value = 18304.666
$ws.Range("J77").Value = 18304.666
$ws.Range("L77").Value = 54913.99800000001
$ws.Range("N77").Value = -63649.99800000001

$ws.Range("H101").Value = 23225
$ws.Range("J101").Value = 23225
$ws.Range("L101").Value = 23225
$ws.Range("N101").Value = -29715

$ws.Range("H110").Value = 38193.332
$ws.Range("J110").Value = 38193.332
$ws.Range("L110").Value = 38193.332
$ws.Range("N110").Value = -46373.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 525.25
$ws.Range("I10").Value = 33.666668
$ws.Range("K10").Value = 101.000004
$ws.Range("M10").Value = 37.999996

$ws.Range("H80").Value = 3596
$ws.Range("J80").Value = 3596
$ws.Range("L80").Value = 10788
$ws.Range("N80").Value = -12660

$ws.Range("H83").Value = 3596
$ws.Range("J83").Value = 3596
$ws.Range("L83").Value = 32364
$ws.Range("N83").Value = -41724

$ws.Range("H129").Value = 2623.5
$ws.Range("I129").Value = 1743.3334
$ws.Range("J129").Value = 2826.6155
$ws.Range("K129").Value = 5230.0002
$ws.Range("L129").Value = 8479.8465
$ws.Range("M129").Value = -230.0002000000004
$ws.Range("N129").Value = -18479.8465

$ws.Range("H131").Value = 922.2432
$ws.Range("I131").Value = 425
$ws.Range("J131").Value = 1160.92
$ws.Range("K131").Value = 1275
$ws.Range("L131").Value = 3482.76
$ws.Range("M131").Value = 3765
$ws.Range("N131").Value = -13562.76

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 12062.5
$ws.Range("J63").Value = 12062.5
$ws.Range("L63").Value = 12062.5
$ws.Range("N63").Value = -13434.5

$ws.Range("H66").Value = 12062.5
$ws.Range("J66").Value = 12062.5
$ws.Range("L66").Value = 36187.5
$ws.Range("N66").Value = -43051.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 5001500
$ws.Range("I12").Value = 10000000
$ws.Range("J12").Value = 3000
$ws.Range("K12").Value = 10000000
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = -9999830
$ws.Range("N12").Value = -3340

$ws.Range("H64").Value = 20150
$ws.Range("J64").Value = 20150
$ws.Range("L64").Value = 20150
$ws.Range("N64").Value = -20600

$ws.Range("H67").Value = 20150
$ws.Range("J67").Value = 20150
$ws.Range("L67").Value = 20150
$ws.Range("N67").Value = -21710

$ws.Range("H108").Value = 14400
$ws.Range("J108").Value = 14400
$ws.Range("L108").Value = 14400
$ws.Range("N108").Value = -22080

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 75253.75
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 100005
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 100005
$ws.Range("M7").Value = -887
$ws.Range("N7").Value = -100231

$ws.Range("H63").Value = 24249.75
$ws.Range("J63").Value = 24249.75
$ws.Range("L63").Value = 24249.75
$ws.Range("N63").Value = -25497.75

$ws.Range("H66").Value = 24249.75
$ws.Range("J66").Value = 24249.75
$ws.Range("L66").Value = 72749.25
$ws.Range("N66").Value = -78989.25

$ws.Range("H76").Value = 13800
$ws.Range("J76").Value = 13800
$ws.Range("L76").Value = 13800
$ws.Range("N76").Value = -14430

$ws.Range("H79").Value = 13800
$ws.Range("J79").Value = 13800
$ws.Range("L79").Value = 13800
$ws.Range("N79").Value = -15984

$ws.Range("H80").Value = 11500
$ws.Range("J80").Value = 11500
$ws.Range("L80").Value = 11500
$ws.Range("N80").Value = -13496

$ws.Range("H83").Value = 11500
$ws.Range("J83").Value = 11500
$ws.Range("L83").Value = 34500
$ws.Range("N83").Value = -44484

$ws.Range("H100").Value = 674.2069
$ws.Range("I100").Value = 566.4545000000001
$ws.Range("J100").Value = 1012.8571
$ws.Range("K100").Value = 1132.909
$ws.Range("L100").Value = 2025.7142
$ws.Range("M100").Value = -591.9090000000001
$ws.Range("N100").Value = -3107.7142
